$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "thermal men leggings",
    "compression pants with knee pads",
    "under armor compression pants men",
    "compression pants with pads",
    "basketball warm up pants men",
    "basketball compression leggings",
    "bjj compression pants",
    "padded leggings basketball",
    "football tights",
    "knee pads for volleyball youth",
    "recovery compression pants men",
    "hockey compression leggings",
    "knee compressions",
    "men sports tights",
    "wrestling knee pads men",
    "mens running leggings",
    "thick knee pad",
    "football pads pants",
    "boys athletic pants",
    "running compression tights",
    "volleyball gear men",
    "girls capri leggings",
    "black capri leggings",
    "coolomg basketball knee pads",
    "tight with knee pads",
    "tights mens nike",
    "muscle compression pants for men",
    "hayabusa compression pants men",
    "girls athletic leggings",
    "black capri leggings for women",
    "asics youth knee pads",
    "men's tights sports",
    "tights with knee pads",
    "lavento compression pants",
    "tsla compression pants men",
    "mens black baseball pants",
    "knee pad protector",
    "compression pants boys basketball",
    "black leggings men",
    "athletic leggings men",
    "youth hex knee pads",
    "wrestling youth knee pads",
    "adult volleyball knee pads",
    "basketball shorts with pads",
    "mens wrestling tights",
    "softball sliding pants",
    "black baseball pants",
    "wrestling tights for men",
    "mens running capris",
    "mens football pants with pads",
    "soccer pants",
    "baseball sliding shorts boys",
    "compression knee guards",
    "knees protector",
    "baseball pants youth",
    "mens basketball shorts black",
    "mens compression",
    "wrestling clothes for men",
    "mens sports pants",
    "boys gym pants",
    "knee pads football adult",
    "leggings men short",
    "running pants men tall",
    "boys running tights youth",
    "softball pants for girls youth",
    "little boys athletic pants",
    "boys running pants size",
    "football youth pants",
    "comfortable knee pads",
    "knee sleeve with pad",
    "women compression tights",
    "knee pads toddler",
    "knee pads basketball kids",
    "knee protector for kids",
    "dark green knee pads for basketball",
    "skateboarding knee pads youth",
    "compression pants with padding basketball",
    "firefighter compression pants",
    "skateboard knee and elbow pads youth",
    "skateboard youth knee pads",
    "warm up pants men",
    "mc david knee pad",
    "compression knee pads men",
    "men basketball pants",
    "premium knee pad",
    "kids compression pants for basketball",
    "toddler knee pad",
    "padded knee sleeves men",
    "mtb knee pads for men",
    "compression knee sleeves pads",
    "basketball knee pads leggings",
    "soccer compression pants men",
    "mens small leggings",
    "compressions pants mens",
    "capri basketball leggings for boys",
    "youth small black baseball pants",
    "mens black leggings",
    "basket ball knee pads youth",
    "bjj leggings men",
    "mens gym tights"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
